$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 32: 2/23/2020 entry
$ws.Range("A32").Value = 43884
$ws.Range("A32").NumberFormat = "m/d/yy"
$ws.Range("B32").Value = "1:00pm - 6:00 pm"
$ws.Range("C32").Value = "Vaishakhi,Anjana"
$ws.Range("D32").Value = "To understand the architecture of the H2 Database from a code design perspective(not necessarily from a Systems design perspective)"
$ws.Range("E32").Value = "Understood both the documented as well as charted the is implemented architecture."
$ws.Range("F32").Value = "It becomes really easy to understand the the system, if it has been documented. In our case, the architecture was documented in the form of text. This was however an older version of the architecture, and worked on charting the is-implemented version. We saw some architectural drift, from what was documented and what is currently being used. I must say that the community of mainters and developers that support H2 database, do a fantastic job at sticking to the architecture, and if there is a change, they systematically report it as well."
$ws.Range("G32").Value = "Excited, as I love architecture design. I wish Andre and Kaj, could give us a lecture on System Design as well. I really want to learn Distributed Computing, and System Design in depth. I would really be excited if we learn about Microservices, system scalability, load balancing, cron scheduling. "

# Row 33: 2/24/2020 entry
$ws.Range("A33").Value = 43885
$ws.Range("A33").NumberFormat = "m/d/yy"
$ws.Range("B33").Value = "9:00 pm - 12:00 pm"
$ws.Range("C33").Value = "Vaishakhi,Anjana"
$ws.Range("D33").Value = "Finish the Social Context part of the homework"
$ws.Range("E33").Value = "Worked on my part of the assignment and collated it with the others' "
$ws.Range("F33").Value = "It was really easy to work on this part, all thanks to the properly maintained H2 application on github. Found how important it is to name issues and pull requests properly, and how important it is to clearly mention fixes in comments. Simple things like spaces,tabs also matter a lot."
$ws.Range("G33").Value = "Excited, as we finished the work pretty early"

# Row heights to match wrapped content
$ws.Rows.Item(32).RowHeight = 255
$ws.Rows.Item(33).RowHeight = 136

# Update view state: scroll + selection
$ws.Application.ActiveWindow.ScrollColumn = 4
$ws.Range("G33").Select()
